$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- New journal entries in rows 22-25 ---

# Row 22: just a date was added (activity/heures stay blank)
$ws.Range("A22").Value = 43195

# Row 23: new entry (two-line activity text -> taller row, like the other ht=30 rows)
$ws.Range("A23").Value = 43200
$ws.Range("B23").Value = "Implémentation de l'historique sur les fonctionnalités outils"
$ws.Range("C23").Value = 13
$ws.Rows.Item(23).RowHeight = 30

# Row 24: new entry
$ws.Range("A24").Value = 43201
$ws.Range("B24").Value = "suite de l'implémentation de l'historique"
$ws.Range("C24").Value = 10

# Row 25: new entry
$ws.Range("A25").Value = 43204
$ws.Range("B25").Value = "finalisation de l'ordre des calques"
$ws.Range("C25").Value = 2

# --- Shift the "Total" row from 34 down to 35 ---

# Preserve the Total row's current content+format by copying it one row down first.
$ws.Range("B34:C34").Copy()
$ws.Range("B35:C35").PasteSpecial(-4122)

# Turn row 34 back into a plain blank data row matching the other blank rows (26-33),
# reusing their formatting so no new cell styles are introduced.
$ws.Range("A26:C26").Copy()
$ws.Range("A34:C34").PasteSpecial(-4122)
$ws.Range("A34").ClearContents()
$ws.Range("B34").ClearContents()
$ws.Range("C34").ClearContents()

# Write the Total label and the updated sum formula on the new row 35.
$ws.Range("B35").Value = "Total"
$ws.Range("C35").Formula = "=SUM(C5:C34)"

# --- Selection as left by the author ---
$ws.Range("K3").Select()
